$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: section header label
$ws.Range("A12").Value = "Rerun, with release mode"

# Row 13: column headers (reuse existing shared strings)
$ws.Range("A13").Value = "String count"
$ws.Range("B13").Value = "Naïve build time (ms)"
$ws.Range("C13").Value = "New merge time (ms)"
$ws.Range("D13").Value = "Interleave iterations"

# Rows 14-21: data
$ws.Range("A14").Value = 2048
$ws.Range("B14").Value = 10.6252
$ws.Range("C14").Value = 14.603300000000001
$ws.Range("D14").Value = 9

$ws.Range("A15").Value = 8192
$ws.Range("B15").Value = 32.121400000000001
$ws.Range("C15").Value = 39.704000000000001
$ws.Range("D15").Value = 10

$ws.Range("A16").Value = 32768
$ws.Range("B16").Value = 149.14869999999999
$ws.Range("C16").Value = 156.85560000000001
$ws.Range("D16").Value = 11

$ws.Range("A17").Value = 131072
$ws.Range("B17").Value = 819.89869999999996
$ws.Range("C17").Value = 732.73979999999995
$ws.Range("D17").Value = 12

$ws.Range("A18").Value = 524288
$ws.Range("B18").Value = 3928.2964000000002
$ws.Range("C18").Value = 3166.4897999999998
$ws.Range("D18").Value = 13

$ws.Range("A19").Value = 2097152
$ws.Range("B19").Value = 20510.476299999998
$ws.Range("C19").Value = 13865.379300000001
$ws.Range("D19").Value = 15

$ws.Range("A20").Value = 4194304
$ws.Range("B20").Value = 53606.484700000001
$ws.Range("C20").Value = 29184.9133
$ws.Range("D20").Value = 14

$ws.Range("A21").Value = 7438776
$ws.Range("B21").Value = 91701.114000000001
$ws.Range("C21").Value = 56427.221700000002
$ws.Range("D21").Value = 15

# Update selection to match final cursor position in the diff (D21)
$ws.Range("D21").Select()
